$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new blank data rows, copying the border/format of an
#    existing "middle" data row into them, so every data row keeps the
#    same medium-box-border look used throughout the table.
#    Final layout needs 12 data rows (3..14) instead of the original 10
#    (3..12), so two rows are added. We insert both just above the last
#    (bottom-bordered) row so the bottom border of the table stays put.
# ---------------------------------------------------------------------------
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).Insert()
$ws.Range("A11:G11").Copy()
$ws.Range("A12:G13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Give columns C and E (and, on the last row, D too) the explicit black
#    font used by the "Good" / "Mid" accuracy columns in the final report.
# ---------------------------------------------------------------------------
$ws.Range("C3:C14").Font.Color = 255
$ws.Range("E3:E14").Font.Color = 255
$ws.Range("D14").Font.Color = 255

# ---------------------------------------------------------------------------
# 3. Write the final header / data values.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Linear Regression"
$ws.Range("C3").Value = -191.81717677875801
$ws.Range("D3").Value = 71.550196435691205
$ws.Range("E3").Value = -172.35457740000001
$ws.Range("F3").Value = -670.23383041305397
$ws.Range("G3").Value = -1383.22188829421

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Linear Regression+PCA"
$ws.Range("C4").Value = 88.715486069492201
$ws.Range("D4").Value = 80.115493957702597
$ws.Range("E4").Value = 84.424648259999998
$ws.Range("F4").Value = 74.006588099760194
$ws.Range("G4").Value = 70.178367723452098

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Lasso Regression"
$ws.Range("C5").Value = 79.420726020000004
$ws.Range("D5").Value = 94.864480171692904
$ws.Range("E5").Value = 79.420726020000004
$ws.Range("F5").Value = 63.765347125414202
$ws.Range("G5").Value = 8.4810702840791308

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Lasso Regression+normalization"
$ws.Range("C6").Value = 69.770022319999995
$ws.Range("D6").Value = 96.131891332527104
$ws.Range("E6").Value = 69.770022319999995
$ws.Range("F6").Value = 77.520144139150702
$ws.Range("G6").Value = 64.940395326255299

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Lasso Regression+normalization+ lag1"
$ws.Range("C7").Value = 69.779661559999994
$ws.Range("D7").Value = 96.136086308924206
$ws.Range("E7").Value = 69.779661559999994
$ws.Range("F7").Value = 77.074776640014207
$ws.Range("G7").Value = 71.002143130205795

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Lasso Regression+normalization+ lag1 +PCA"
$ws.Range("C8").Value = 75.625968420000007
$ws.Range("D8").Value = 90.313239586760503
$ws.Range("E8").Value = 75.625968420000007
$ws.Range("F8").Value = 84.203732293872307
$ws.Range("G8").Value = 83.519033481304604

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Polynomial Regression"
$ws.Range("C9").Value = -77.956117480000003
$ws.Range("D9").Value = 69.781725315391697
$ws.Range("E9").Value = -77.956117480000003
$ws.Range("F9").Value = -273.85325356476699
$ws.Range("G9").Value = -390.31367793150099

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Polynomial Regression+normalization"
$ws.Range("C10").Value = 66.510769589999995
$ws.Range("D10").Value = 89.384114536514403
$ws.Range("E10").Value = 66.510769589999995
$ws.Range("F10").Value = 77.520144139150702
$ws.Range("G10").Value = 64.940395326255299

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Random Forest"
$ws.Range("C11").Value = 88.329189929999998
$ws.Range("D11").Value = 95.496711032342702
$ws.Range("E11").Value = 88.329189929999998
$ws.Range("F11").Value = 74.736648783412605
$ws.Range("G11").Value = 59.5861883653508

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Decision Tree Regression"
$ws.Range("C12").Value = 81.049431760000004
$ws.Range("D12").Value = 95.4664193060201
$ws.Range("E12").Value = 81.049431760000004
$ws.Range("F12").Value = 49.7656108508356
$ws.Range("G12").Value = 25.156157702698899

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Decision Tree Regression+feature selection"
$ws.Range("C13").Value = 81.739564189999996
$ws.Range("D13").Value = 96.116233891006303
$ws.Range("E13").Value = 81.739564189999996
$ws.Range("F13").Value = 85.461516214202803
$ws.Range("G13").Value = 78.441738356474701

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "GB Regression"
$ws.Range("C14").Value = 86.770905940000006
$ws.Range("D14").Value = 97.665839938404105
$ws.Range("E14").Value = 86.770905940000006
$ws.Range("F14").Value = 62.5372804024381
$ws.Range("G14").Value = 51.183324516805101

# ---------------------------------------------------------------------------
# 4. Add the summary AVERAGE row.
# ---------------------------------------------------------------------------
$ws.Range("C15").Formula = "=AVERAGE(C4:C14)"
$ws.Range("D15:G15").Formula = "=AVERAGE(D3:D14)"

# ---------------------------------------------------------------------------
# 5. Column width tweaks (column B widened, column F narrowed slightly).
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 44.2857142857143
$ws.Columns.Item(6).ColumnWidth = 16.8571428571429

# ---------------------------------------------------------------------------
# 6. Selection moves to C18 in the saved file.
# ---------------------------------------------------------------------------
$ws.Range("C18").Select()
